$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.387351989746094
$ws.Range("B1").Value = 1.964551687240601
$ws.Range("C1").Value = 2.871350288391113
$ws.Range("D1").Value = 4.797849655151367
$ws.Range("E1").Value = 0.9758411645889282
